# repull data, push all data, mean calculation
# Update column F (dSF) values for the affected rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = 3
$ws.Range("F3").Value  = -2
$ws.Range("F4").Value  = 3
$ws.Range("F5").Value  = 8
$ws.Range("F6").Value  = -3
$ws.Range("F7").Value  = 6
$ws.Range("F8").Value  = 3
$ws.Range("F9").Value  = -1
$ws.Range("F12").Value = -3
$ws.Range("F14").Value = -1
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = 10
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = -5
